# Generated by Katalon AI - apply renaming + column width adjustments

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (columns E..L / 5..12) ---
# Note: ColumnWidth in character units maps to raw OOXML width with a
# +5/6 (0.8333...) offset added by the runtime, so subtract that offset
# to land exactly on the desired stored width.
$offset = 0.8333333333333334

$ws.Columns.Item(5).ColumnWidth  = 32 - $offset
$ws.Columns.Item(6).ColumnWidth  = 34 - $offset
$ws.Columns.Item(7).ColumnWidth  = 34 - $offset
$ws.Columns.Item(8).ColumnWidth  = 34 - $offset
$ws.Columns.Item(10).ColumnWidth = 44 - $offset
$ws.Columns.Item(11).ColumnWidth = 34 - $offset
$ws.Columns.Item(12).ColumnWidth = 31 - $offset

# --- Header text (row 1) renames ---
$ws.Range("A1").Value = "button_closeAlert_class"
$ws.Range("B1").Value = "button_closeAlert_class_1"
$ws.Range("C1").Value = "div_backdropComponents_class"
$ws.Range("D1").Value = "div_backdropComponents_class_1"
$ws.Range("E1").Value = "div_testSuiteItem_internalText"
$ws.Range("F1").Value = "div_testSuiteItem_internalText_1"
$ws.Range("G1").Value = "div_testSuiteItem_internalText_2"
$ws.Range("H1").Value = "div_testSuiteItem_internalText_3"
$ws.Range("J1").Value = "link_organizationLink_internalRoleLinkName"
$ws.Range("K1").Value = "link_organizationLink_project_id"
$ws.Range("L1").Value = "link_organizationLink_team_id"
